$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Ensure the new rows are written as plain text (avoid Excel auto-converting
# numeric-looking statute codes like "4510.11" or currency-looking "$ 0" strings)
$ws.Range("A660:K678").NumberFormat = "@"

$newRows = @(
    @('21TRD09437', 'Hemmeter', 'DUS', '4510.11', 'M1', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Hemmeter', '1ST SPEED 1 YR SCHOOL >35MPHM4', '4511.21B1A', 'M4', 'Dismissed', $null, ' ', ' ', $null, $null),
    @('21TRD09437', 'Hemmeter', 'RECKLESS OPERATION 1ST IN 1 YR', '4511.20', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Hemmeter', 'DUS', '4510.11', 'M1', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Hemmeter', '1ST SPEED 1 YR SCHOOL >35MPHM4', '4511.21B1A', 'M4', 'Dismissed', $null, ' ', ' ', $null, $null),
    @('21TRD09437', 'Hemmeter', 'RECKLESS OPERATION 1ST IN 1 YR', '4511.20', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', 'DUS', '4510.11', 'M1', 'No Contest', 'Not Guilty - Allied Offense', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', '1ST SPEED 1 YR SCHOOL >35MPHM4', '4511.21B1A', 'M4', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', 'RECKLESS OPERATION 1ST IN 1 YR', '4511.20', 'MM', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', 'DUS', '4510.11', 'M1', 'No Contest', 'Guilty', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', '1ST SPEED 1 YR SCHOOL >35MPHM4', '4511.21B1A', 'M4', 'No Contest', 'Guilty - Allied Offense', '$ 0', '$ 0', $null, $null),
    @('21TRD09437', 'Bunner', 'RECKLESS OPERATION 1ST IN 1 YR', '4511.20', 'MM', 'Dismissed', $null, ' ', ' ', $null, $null),
    @('21TRD09437', 'Bunner', 'DUS', '4510.11', 'M1', 'Guilty', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('21TRD09437', 'Bunner', '1ST SPEED 1 YR SCHOOL >35MPHM4', '4511.21B1A', 'M4', 'Guilty', 'Guilty - Allied Offense', '$ 0', '$ 0', 'None', 'None'),
    @('21TRD09437', 'Bunner', 'RECKLESS OPERATION 1ST IN 1 YR', '4511.20', 'MM', 'Guilty', 'Guilty', '$ 0', '$ 0', 'None', 'None'),
    @('22CRB00098', 'Bunner', 'DOMESTIC VIOLENCE  M1', '2919.25A', 'M1', 'Not Guilty', $null, $null, $null, $null, $null),
    @('22CRB00098', 'Bunner', 'ASSAULT', '2903.13', 'M1', 'Not Guilty', $null, $null, $null, $null, $null),
    @('22CRB00098', 'Bunner', 'DOMESTIC VIOLENCE  M1', '2919.25A', 'M1', 'Guilty', $null, $null, $null, $null, $null),
    @('22CRB00098', 'Bunner', 'ASSAULT', '2903.13', 'M1', 'Guilty', $null, $null, $null, $null, $null),
)

$startRow = 660
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $val = $rowData[$c]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}
